$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the template placeholder text that was wrongly replaced with literal
# values instead of the ${PLACEHOLDER_n} tokens.
$ws.Range("B1").Value = '${PLACEHOLDER_3}'
$ws.Range("J10").Value = 'It''s a ${PLACEHOLDER_1}'
$ws.Range("D27").Value = 'It''s a ${PLACEHOLDER_2}'

# Restore the selected cell to B2 (was incorrectly left at C1).
$ws.Range("B2").Select()
